$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing the existing data (old rows 2-16)
# down to rows 3-17.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "Juist, Inselgemeinde" record.
$ws.Range("A2").Value = 452013
$ws.Range("B2").Value = "Juist, Inselgemeinde"
$ws.Range("C2").Value = 53.678347
$ws.Range("D2").Value = 6.995328
$ws.Range("E2").Value = "https://drive.google.com/drive/folders/0BxMfdWAA8UdsUFhRMzdDdnh4Z1E"
$ws.Range("F2").Value = "http://www.gemeinde-juist.de/"

# Match the author's row height for the new row.
$ws.Rows.Item(2).RowHeight = 13.8

# Update the selection to highlight the newly inserted row, matching the
# author's recorded UI state after the edit.
$ws.Rows.Item(2).Select()

# The sheet's filter-database defined name tracks the used range and grows
# by one row along with the data.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "BLP-URLs!_FilterDatabase") {
        $n.RefersTo = "='BLP-URLs'!`$A`$1:`$G`$405"
    }
}
